$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario")

# Header text fix
$ws.Range("B1").Value = "Non-Residential"

# Updated projection values (new parquet files structure/values)
$ws.Range("B2").Value = 5989620461.822261
$ws.Range("C2").Value = 15539725146.071728

$ws.Range("B3").Value = 6189508557.828688
$ws.Range("C3").Value = 15860659116.453512

$ws.Range("B4").Value = 5731297592.470088
$ws.Range("C4").Value = 15469971370.55126

$ws.Range("B5").Value = 5604720706.924748
$ws.Range("C5").Value = 15512862022.853907

$ws.Range("B6").Value = 5374894548.775076
$ws.Range("C6").Value = 15360897580.2493

$ws.Range("B7").Value = 5121955169.398381
$ws.Range("C7").Value = 15157003643.782705

$ws.Range("B8").Value = 4853096961.886391
$ws.Range("C8").Value = 14928953612.305525
